# Append a new daily price row (2025-04-14) to every price sheet in the
# workbook. Each sheet's latest entry (2025-04-13, row 43) is carried
# forward unchanged into the new row 44 - only the date in column A
# actually changes; the price in column B is identical to the prior day
# on every sheet.
#
# Row 43 -> new row 44 price carried forward, per sheet (workbook order):
#   N-Dense                   : 40
#   N-Type                    : 41.5
#   N-type Wafer               : 1.25
#   Cell Topcon 183mm          : 0.303
#   Module Topcon 183mm        : 0.1
#   Silver Rear_side           : 5,192
#   Silver Busbar front-side   : 7,773
#   Silver finger front-side   : 7,823
#   USD_CNY                    : 7.3258

$wb = $excel.ActiveWorkbook

$newDateFormula = "=""2025-04-14"""
$lastRow = 43
$newRow = 44
$scratchRow = 1000
$scratchCol = 100

foreach ($ws in $wb.Worksheets) {
    # Duplicate the last data row (date + price) down into the new row.
    # A plain Copy preserves the existing cell content/type exactly (both
    # cells stay plain text, same as every other row in these columns)
    # with no incidental formatting/style changes - and since the price
    # is unchanged day over day on every sheet, column B is already
    # correct after this.
    $src = $ws.Range("A$lastRow`:B$lastRow")
    $dst = $ws.Range("A$newRow`:B$newRow")
    $src.Copy($dst)

    # Column A still needs the new date. Writing the date string straight
    # into the cell's Value would get auto-coerced into a date serial by
    # Excel's usual parsing, unlike the existing (text) date cells. Instead,
    # compute it as a text formula result off-sheet (a formula's string
    # result is never re-interpreted as a date/number) and paste just the
    # value in, which keeps the cell a plain text cell like its neighbours.
    $scratch = $ws.Cells.Item($scratchRow, $scratchCol)
    $scratch.Formula = $newDateFormula
    $scratch.Copy()
    $ws.Cells.Item($newRow, 1).PasteSpecial(-4163)  # xlPasteValues
    $scratch.Clear()
    $excel.CutCopyMode = $false
}
